$wb = $excel.ActiveWorkbook

# --- Rushing sheet updates ---
$rushing = $wb.Worksheets.Item("Rushing")

# K.Murray (row 2)
$rushing.Range("C2").Value = 26
$rushing.Range("D2").Value = 22
$rushing.Range("E2").Value = 13

# J.Conner (row 5)
$rushing.Range("C5").Value = 86
$rushing.Range("D5").Value = 45
$rushing.Range("E5").Value = 24
$rushing.Range("F5").Value = 34

# C.Kirk (row 8)
$rushing.Range("C8").Value = 6

# R.Moore (row 9)
$rushing.Range("C9").Value = 11

# --- Receiving sheet updates ---
$receiving = $wb.Worksheets.Item("Receiving")

# C.Edmonds (row 2) - zeroed out
$receiving.Range("C2").Value = 0
$receiving.Range("D2").Value = 0
$receiving.Range("E2").Value = 0
$receiving.Range("F2").Value = 0
$receiving.Range("G2").Value = 0
$receiving.Range("H2").Value = 0

# J.Conner (row 3)
$receiving.Range("C3").Value = 24
$receiving.Range("D3").Value = 22
$receiving.Range("G3").Value = 2
$receiving.Range("H3").Value = 1

# D.Hopkins (row 5) - zeroed out
$receiving.Range("C5").Value = 0
$receiving.Range("D5").Value = 0
$receiving.Range("E5").Value = 0
$receiving.Range("F5").Value = 0
$receiving.Range("G5").Value = 0
$receiving.Range("H5").Value = 0

# A.Green (row 6)
$receiving.Range("C6").Value = 48
$receiving.Range("D6").Value = 30
$receiving.Range("E6").Value = 19
$receiving.Range("F6").Value = 13
$receiving.Range("G6").Value = 12

# C.Kirk (row 7)
$receiving.Range("C7").Value = 45
$receiving.Range("D7").Value = 36
$receiving.Range("E7").Value = 19
$receiving.Range("F7").Value = 11

# R.Moore (row 8)
$receiving.Range("C8").Value = 51
$receiving.Range("D8").Value = 43
$receiving.Range("E8").Value = 5

# D.Harris (row 13)
$receiving.Range("C13").Value = 44
$receiving.Range("D13").Value = 37
$receiving.Range("G13").Value = 7
$receiving.Range("H13").Value = 4

# Make the Receiving sheet the active sheet / tab, with H3 selected
$receiving.Activate()
$receiving.Range("H3").Select()
